$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume/1h (E) values for rows 2-49
$ws.Range("D2").Value = "27.967.55"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "1.762.60"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.10"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3773"
$ws.Range("E7").Value = "  -4.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3354"
$ws.Range("E8").Value = "  -4.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.80"
$ws.Range("E9").Value = "  -4.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.123"
$ws.Range("E10").Value = "  -6.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07173"
$ws.Range("E11").Value = "  -5.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.38"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.189"
$ws.Range("E14").Value = "  -5.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.173"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "1.757.65"
$ws.Range("E16").Value = "  -3.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001053"
$ws.Range("E17").Value = "  -5.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06581"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.43"
$ws.Range("E19").Value = "  -6.09%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.94"
$ws.Range("E21").Value = "  -5.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.272"
$ws.Range("E22").Value = "  -5.01%  "
$ws.Range("D23").Value = "27.953.73"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.71"
$ws.Range("E24").Value = "  -8.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.369"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.38"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.83"
$ws.Range("E27").Value = "  -7.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.332"
$ws.Range("E28").Value = "  -9.93%  "
$ws.Range("D29").Value = "1.957.92"
$ws.Range("E29").Value = "  -3.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.264"
$ws.Range("E30").Value = "  -15.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.94"
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.021"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.788"
$ws.Range("E33").Value = "  -7.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08749"
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.22"
$ws.Range("E35").Value = "  -8.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02337"
$ws.Range("E36").Value = "  -4.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6569"
$ws.Range("E37").Value = "  -6.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06191"
$ws.Range("E38").Value = "  -6.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.148"
$ws.Range("E39").Value = "  -7.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2106"
$ws.Range("E40").Value = "  -5.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.211"
$ws.Range("E41").Value = "  -4.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.452"
$ws.Range("E42").Value = "  -10.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.021"
$ws.Range("E43").Value = "  -6.16%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.69"
$ws.Range("E45").Value = "  -6.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.830"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6034"
$ws.Range("E47").Value = "  -7.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.94"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.011"
$ws.Range("E49").Value = "  -7.69%  "

# Rows 50 and 51: coins swap order (Cronos <-> EOS) plus updated values
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.183"
$ws.Range("E50").Value = "  +2.02%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07161"
$ws.Range("E51").Value = "  -0.91%  "
